$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("B$row").Value = 3.230985683306322
    $ws.Range("C$row").Value = 1.667794583268128
    $ws.Range("D$row").Value = 0.8054896365839992
    $ws.Range("E$row").Value = 0.496779210170732
    $ws.Range("G$row").Value = 6.201049113329182
}
